$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bugs")
$lo = $ws.ListObjects.Item("Table1")

# ------------------------------------------------------------------
# 1. Grow the table to A1:F8 (adds column F and rows up to 8)
# ------------------------------------------------------------------
$lo.Resize($ws.Range("A1:F8"))

# ------------------------------------------------------------------
# 2. Column widths: insert a "User" column (D) between Device(C) and
#    the old Notes column (now E); Status becomes F.
# ------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 54.7109375
$ws.Columns("D").ColumnWidth = 35.28515625
$ws.Columns("E").ColumnWidth = 54.7109375
$ws.Columns("F").ColumnWidth = 26.140625

# ------------------------------------------------------------------
# 3. Header row (row 1) - shift Notes/Status right, add "User"
# ------------------------------------------------------------------
$ws.Range("D1").Value = "User"
$ws.Range("E1").Value = "Notes"
$ws.Range("F1").Value = "Status"

# ------------------------------------------------------------------
# 4. Data rows
# ------------------------------------------------------------------
# Row 2 (bug #1)
$ws.Range("E2").Value = "Should be the same as bug #5"
$ws.Range("F2").Value = "OK"

# Row 3 (bug #2)
$ws.Range("D3").Value = "Elena"
$ws.Range("E3").Value = "The issue is located on the HTML builder, the property width is missing"
$ws.Range("F3").Value = "OK"

# Row 4 (bug #3)
$ws.Range("D4").Value = "Lembo"
$ws.Range("F4").Value = "OK"

# Row 5 (bug #4)
$ws.Range("D5").Value = "Lembo"
$ws.Range("E5").Value = "Can't reproduce"

# Row 6 (bug #5)
$ws.Range("C6").Value = "Android > 13"
$ws.Range("D6").Value = "Fede, Mauro"
$ws.Range("E6").Value = "With Android > 13, it's not possible anymore to use READ_EXTERNAL_STORAGE, but READ_MEDIA_IMAGES and other two permissions must be used"
$ws.Range("F6").Value = "OK"

# Row 7 (bug #6)
$ws.Range("D7").Value = "Lembo"
$ws.Range("F7").Value = "OK"

# Row 8 (bug #7) - brand new row
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Alla creazione della notifica avviene un crash"
$ws.Range("C8").Value = "Galaxy S22"
$ws.Range("D8").Value = "Fede, Mauro"
$ws.Range("E8").Value = 'With Android > 13, it is necessary to add the permissions <uses-permission android:name="android.permission.USE_EXACT_ALARM" />'
$ws.Range("F8").Value = "OK"

# ------------------------------------------------------------------
# 5. Apply the "Good" (green) row style to every row whose Status is
#    OK, by copying the formatting of an already-green row (row 3).
# ------------------------------------------------------------------
$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A2:F2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A4:F4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A6:F6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A7:F7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A8:F8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A3:F3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 6. Data validation dropdown on Status column moves from E to F
# ------------------------------------------------------------------
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("F2:F1048576").Validation.Add(3, 1, 1, "OK,KO,WIP,ON HOLD") | Out-Null
$ws.Range("E2:E1048576").Validation.Delete()

# ------------------------------------------------------------------
# 7. Selection / active cell, to match the authored state
# ------------------------------------------------------------------
$ws.Range("A7:F7").Select() | Out-Null
$ws.Application.ActiveWindow.RangeSelection
$ws.Activate()
